$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Absent" column (H) with the consolidated values.
# Absent = 1 - Real (column E)
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
